# One positive test for Login page:
# add a new row for "LoginWithNonExistingUserButValidPassword" with a valid
# (but non-existing) email and an "Existing" password, mirroring the style
# and hyperlink pattern used by the other Email cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "LoginWithNonExistingUserButValidPassword"
$ws.Range("C6").Value = "abv@abv.bg"
$ws.Range("D6").Value = "Existing"

# Turn the email cell into a mailto hyperlink, then give it the same
# "Hyperlink" cell style used by the other Email column entries.
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:abv@abv.bg")
$ws.Range("C6").Style = $ws.Range("C5").Style

$ws.Range("D6").Select()
